# Reorder the comma-separated "Recorded By" names in column G:
#   - if the literal (case-sensitive) token "System" is present, move it to
#     the front and reverse the order of the remaining tokens
#   - otherwise, simply reverse the order of all tokens
# Single-token values are left untouched (reversal is a no-op for them).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

function Reverse-Array($arr) {
    $result = @()
    for ($i = $arr.Length - 1; $i -ge 0; $i--) {
        $result += $arr[$i]
    }
    return $result
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -eq "") { continue }

    $parts = $val -split ','
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    if ($parts.Length -lt 2) { continue }

    $systemIndex = -1
    for ($i = 0; $i -lt $parts.Length; $i++) {
        if ($parts[$i].Equals('System')) {
            $systemIndex = $i
            break
        }
    }

    if ($systemIndex -ge 0) {
        $rest = @()
        for ($i = 0; $i -lt $parts.Length; $i++) {
            if ($i -ne $systemIndex) { $rest += $parts[$i] }
        }
        $rest = Reverse-Array $rest
        $newParts = @('System') + $rest
    } else {
        $newParts = Reverse-Array $parts
    }

    $newVal = [string]::Join(', ', $newParts)

    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}
